# Apply the crypto price/volume refresh captured in the commit diff.
# D-column price cells are free-form text (e.g. "29.776.07", "8.500",
# "12.50") rather than numbers, so writing through .Value while the
# cell is text-formatted ("@") preserves the exact digits/trailing
# zeros; Style is then reset back to "Normal" so no stray cell style
# index is left behind (matches the original, unstyled data cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.776.07'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +7.02%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.948.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.39%  '

$ws.Range('E4').Value = '  -0.65%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '341.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.78%  '

$ws.Range('E6').Value = '  -0.56%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4785'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.81%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4130'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.72%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.80'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.14%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08229'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.95%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.036'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.80%  '

$ws.Range('E12').Value = '  +6.65%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.959.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.70%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.151'
$ws.Range('D14').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.376'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.03%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001057'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.62%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06679'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E20').Value = '  +4.14%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '29.724.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.88%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.580'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.25%  '

$ws.Range('E24').Value = '  +3.68%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.295'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.20%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.185.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.21%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '161.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.71%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.171'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.85%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.644'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.94%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '123.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.55%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.009'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.63%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09649'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.31%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.473'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.86%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.687'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.498'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.37%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06261'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.90%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02316'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.54%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.500'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.14%  '

$ws.Range('E40').Value = '  +2.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6074'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.40%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.21%  '

$ws.Range('E43').Value = '  -0.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1893'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.43%  '

$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.274'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.66%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.390'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +32.33%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5717'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.75%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07418'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.989'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.60%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.02'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.03%  '
